$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the existing row 533 (shifts rows 533:653 down to 534:654)
$ws.Rows("533:533").Insert()

# Populate the newly inserted row with the new weekly price record
$ws.Cells.Item(533, 1).Value = 10
$ws.Cells.Item(533, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(533, 3).Value = "La Araucanía"
$ws.Cells.Item(533, 4).Value = 45211
$ws.Cells.Item(533, 5).Value = 9
$ws.Cells.Item(533, 6).Value = "Fruta"
$ws.Cells.Item(533, 7).Value = 100108
$ws.Cells.Item(533, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(533, 9).Value = 100108002
$ws.Cells.Item(533, 10).Value = "Mango"
$ws.Cells.Item(533, 11).Value = "Sin especificar"
$ws.Cells.Item(533, 12).Value = "Primera"
$ws.Cells.Item(533, 13).Value = 1050
$ws.Cells.Item(533, 14).Value = 10000
$ws.Cells.Item(533, 15).Value = 10000
$ws.Cells.Item(533, 16).Value = 10000
$ws.Cells.Item(533, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(533, 18).Value = "Brasil"
$ws.Cells.Item(533, 19).Value = 2500
$ws.Cells.Item(533, 20).Value = 4
